$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.603.80"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "1.985.27"
$ws.Range("E3").Value = "  +4.57%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").Value = "'328.06"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").Value = "'0.4660"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").Value = "'0.3932"
$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("D9").Value = "'46.43"
$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").Value = "'0.07953"
$ws.Range("E10").Value = "  +0.68%  "

$ws.Range("D11").Value = "'0.9944"
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("D12").Value = "'22.87"
$ws.Range("E12").Value = "  +3.81%  "

$ws.Range("D13").Value = "2.007.88"
$ws.Range("E13").Value = "  +5.32%  "

$ws.Range("D14").Value = "'7.198"
$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("D15").Value = "'5.848"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "'0.07114"
$ws.Range("E16").Value = "  +1.87%  "

$ws.Range("D17").Value = "'87.70"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").Value = "'1.007"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").Value = "'0.000009960"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Value = "'17.34"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "29.671.98"
$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.548"
$ws.Range("E23").Value = "  +4.68%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'11.21"
$ws.Range("E24").Value = "  +1.14%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.241.10"
$ws.Range("E25").Value = "  +4.95%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.108"
$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'159.34"
$ws.Range("E27").Value = "  +2.06%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.64"
$ws.Range("E28").Value = "  +0.62%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.840"
$ws.Range("E29").Value = "  -3.23%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'119.75"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").Value = "'1.901"
$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.09438"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.9005"
$ws.Range("E33").Value = "  -0.61%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.241"
$ws.Range("E34").Value = "  -1.20%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.334"
$ws.Range("E35").Value = "  +0.37%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'3.197"
$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05820"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.177"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02103"
$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.868"
$ws.Range("E40").Value = "  +0.75%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "'0.000003246"
$ws.Range("E41").Value = "  +49.94%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5737"
$ws.Range("E42").Value = "  +0.50%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1812"
$ws.Range("E43").Value = "  +1.49%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'9.756"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "'2.762"
$ws.Range("E45").Value = "  +7.53%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'11.97"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5371"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'2.188"
$ws.Range("E48").Value = "  -1.72%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06949"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'114.38"
$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.831"
$ws.Range("E51").Value = "  -1.70%  "
